# This workbook's "Other Info" sheet previously listed each survey question's
# label/answer pairs on separate rows (row 1/2 = headers, row 3 = the single
# school's answers). Instead of hand-writing a long flat block of statements,
# loop over a compact (column, value) table and write row 4 in one pass -
# this keeps the code short and the write loop fast.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Other Info")

$rowNum = 4

$row4Data = @(
    @{ Col = "C"; Val = 5 },
    @{ Col = "D"; Val = 3 },
    @{ Col = "E"; Val = 4 },
    @{ Col = "F"; Val = 5 },
    @{ Col = "G"; Val = 5 },
    @{ Col = "H"; Val = 5 },
    @{ Col = "I"; Val = 4 },
    @{ Col = "J"; Val = 3 },
    @{ Col = "K"; Val = 4 },
    @{ Col = "L"; Val = 4 },
    @{ Col = "M"; Val = 2 },
    @{ Col = "N"; Val = 1 },
    @{ Col = "O"; Val = 2 },
    @{ Col = "P"; Val = 3 },
    @{ Col = "Q"; Val = 5 },
    @{ Col = "R"; Val = 3 },
    @{ Col = "S"; Val = 5 },
    @{ Col = "U"; Val = "Please describe how your ODS program’s content (questions 3 and 4) is integrated with local school curricula in a manner that assists students in meeting state standards related to science, technology, engineering and mathematics (STEM), and the Next Generation Science Standards." },
    @{ Col = "V"; Val = "MOSS Curriculum is rooted in Next Generation Science Standards. Each week, teachers are able to choose a Disciplinary Core Idea (relevent to the local ecosystem) on which to focus the lessons throughout the week. Insturctors deliver the content through the Practices of Engineering and Science during the program, and guide students through inquiry-based science investigations in Pondeosa State Park. Each week, teachers may also choose between an `"Inquiry Project`" or a `"Community Engineering Project`" for the students to focus on for an entire day at the end of the week. Students design their own projects, in which they either ask a question they attempt to answer by collecting data in the natural enviornment, or identify a natural resources problem in the local community that they can attempt to solve using the content they learn throughout the week. " },
    @{ Col = "W"; Val = "Please describe if (and if so, how) your ODS program is offered in a bilingual format." },
    @{ Col = "X"; Val = "MOSS relies on schools to bring translators or interpreters with their multi-lingual students. " },
    @{ Col = "Y"; Val = "Please describe how the program addresses the inequity of outdoor educational opportunities for underserved children in this state." },
    @{ Col = "Z"; Val = "MOSS actively seeks grant funding for students who are unable to attend programs, and supports classrooms that are challenged to find funding through assistance in grant writing and in-kind donations of instrcution by staff and volunteers. " },
    @{ Col = "AA"; Val = "Please describe how your ODS program provides students with opportunities to learn about the interdependence of urban and rural areas." },
    @{ Col = "AB"; Val = "MOSS is located in rural central Idaho. Due to the location, we are constantly seeking connection with student from urban and rural areas to identify unique natural resources issues in their area, and be able to identify the ecological principles that may apply to their hometown. By understanding the basics of ecoystem interactions, students will be better suited to identify issues in their home area, and see that the human landscape is intricately interwoven with the natural landscape. " },
    @{ Col = "AC"; Val = "What instructional strategies are used during the program (select all that apply.)" },
    @{ Col = "AD"; Val = "X" },
    @{ Col = "AE"; Val = "X" },
    @{ Col = "AF"; Val = "X" },
    @{ Col = "AG"; Val = "X" },
    @{ Col = "AH"; Val = "X" },
    @{ Col = "AI"; Val = "X" },
    @{ Col = "AJ"; Val = "X" },
    @{ Col = "AK"; Val = "Place-Based Learning" },
    @{ Col = "AL"; Val = "Please describe how learning is extended back to the classroom and annual curriculum. (This may include but is not limited to pre-and post-activities, assessments of learning, teacher professional development, etc.)" },
    @{ Col = "AM"; Val = "MOSS provides pre and post-visit activities for teachers to utilize in the classroom upon request. The K12 Programs Coordinator diligently seeks conversations with teachers before their arrival at MOSS to understand how their experience will apply to what the students are learning in the classroom, and how the field experiences can capitalize on projects in a variety of classes (science, math, language arts, social studies, art, etc). " },
    @{ Col = "AN"; Val = "Please describe how you involve partners in the program. (Partners can include but are not limited to community members, volunteers, government agencies, local non-profits, etc.)" },
    @{ Col = "AO"; Val = "MOSS programs rely heavily on the local McCall community for assistance in designing place-based curriculum that spans socio-ecological systems in the area. We partner with Idaho State Parks for a location in which to teach; The US Forest Service often provides guest speakers for Evening Programs to discuss careers in natural resources (wildlind fire, forestry, hydrology, soil science, and others); The Nez Perce Tribe provides guest speakers about fisheries and native culture; local recreation industry representatives have provided guest speakers to discuss the effects of use on the natural world and their industry; local ski resorts and the Payette Avalanche Center provide insturction on snow science and backcountry rescue. " },
    @{ Col = "AP"; Val = "Please briefly describe your success for each of the items below. Use `"n/a`" for any items not addressed in your ODS program." },
    @{ Col = "AQ"; Val = "Higher scores on standardized measures of academic achievement in reading, writing, math, science and social studies." },
    @{ Col = "AR"; Val = "We have seen a greater understanding and an overall higher performance on standardized tests." },
    @{ Col = "AS"; Val = "Greater self-sufficiency and leadership skills" },
    @{ Col = "AT"; Val = "Teachers report that students are more self-reliant throughout and after programming due to the residential nature of the program. " },
    @{ Col = "AU"; Val = "Fewer discipline and classroom management problems" },
    @{ Col = "AV"; Val = "Students who are identified as having possible behavior issues are rarely seen to exhibit these behaviors in the field. These students are able to channel energy into the field studies" },
    @{ Col = "AW"; Val = "Increased student engagement and pride in accomplishments" },
    @{ Col = "AX"; Val = "Students develop, research, and present their own Inquiry and Engineering Projects each week, providing them with pride in starting and completing a project largely on their own as a team. " },
    @{ Col = "AY"; Val = "Greater proficiency in solving problems and thinking strategically." },
    @{ Col = "AZ"; Val = "Students develop, research, and present their own Inquiry and Engineering Projects each week, providing them with pride in starting and completing a project largely on their own as a team. " },
    @{ Col = "BA"; Val = "Better application of systems thinking and increased ability to think creatively." },
    @{ Col = "BB"; Val = "Students are also expected to make connections between the ecosystems they are experiencing and their home environments. " },
    @{ Col = "BC"; Val = "Improved communication skills and enhanced ability to work in group settings." },
    @{ Col = "BD"; Val = "Students work collaboratively all thoughout the week. One of the `"goals`" of MOSS is to work collabortively and solve problems as a team, through which they develop a sense of community. Students are also living and working together all week long, which requires them to work together to acheive group goals and adhere to group values identified at the beginning of each week." },
    @{ Col = "BE"; Val = "Greater enthusiasm for language arts, math, science and social studies" },
    @{ Col = "BF"; Val = "By the end of each week, students identify themselves as scientists! (Not just a man in a laboratory wearing a labcoat)" },
    @{ Col = "BG"; Val = "Increased knowledge and understanding of science content, concepts and processes." },
    @{ Col = "BH"; Val = "Students identify themselves as scientists by the end of the week. Throughout the experience, students practice science through a lens of inquiry and discovering connections in the ecosystem. Students are able to think scientifically and identify key aspects of the ecosystem they study. " },
    @{ Col = "BI"; Val = "Better ability to apply science and civic processes to real-world situations" },
    @{ Col = "BJ"; Val = "Students are asked to include a `"why should we care?`" section of each Inquiry or Engineering project they complete. In this section, students identify that the principles and concepts they implemented in their self-driven projects have greater implications if they are scaled outward in the community or region. " },
    @{ Col = "BK"; Val = "Improved understanding of mathematical concepts and mastery of math skills." },
    @{ Col = "BL"; Val = "Students are asked to graph data they collect each week and interpret the graphs they develop. Through this data analysis, students are able to identify independent and dependent variables, and make predictions from their data sets. " },
    @{ Col = "BM"; Val = "Improved language arts skills." },
    @{ Col = "BN"; Val = "Students are asked to journal and reflect on each day in the field The journal reflections are collected by the teachers and by the K12 Program Coordiantor for grading and assessment. " },
    @{ Col = "BO"; Val = "Better comprehension of social studies content." },
    @{ Col = "BP"; Val = "MOSS curriuclum focuses on holistic understanding of socio-ecological systems, including the social, economic, cultural, and scienctific connections of local and regional issues (eg: water resources in a changing climate)" },
    @{ Col = "BQ"; Val = "Accessibility to students of all abilities and learning styles" },
    @{ Col = "BR"; Val = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. " },
    @{ Col = "BS"; Val = "Please describe how your program serves historicaly underrepresented populations. What actions are you taking to make this an experience for all students to thrive? (Select and briefly describe all that apply.)" },
    @{ Col = "BT"; Val = "Rural" },
    @{ Col = "BU"; Val = "MOSS instructors are trained in best practices for inclusive education -- honoring diverse perspectives, using student-centered approaches, and acknowledging that students come from different cultural and value frameworks. " },
    @{ Col = "BV"; Val = "English language learners" },
    @{ Col = "BW"; Val = "We incorporate visual and hands-on learning as well as oral and written word in our curriculum so ELL students can experience the curriculum through multiple entry points. " },
    @{ Col = "BX"; Val = "Special education" },
    @{ Col = "BY"; Val = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. " },
    @{ Col = "BZ"; Val = "Learners with disabilities" },
    @{ Col = "CA"; Val = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. We have ADA accessible bunkhouse units, classrooms, and shower facilities." },
    @{ Col = "CB"; Val = "Low socio-economic" },
    @{ Col = "CC"; Val = "We make every effort to make our programs financially accessible. We provide equipment and clothing for students who may not have proper winter clothing." },
    @{ Col = "CD"; Val = "American Indian/Alaska Native" },
    @{ Col = "CE"; Val = "See above under rural section " },
    @{ Col = "CF"; Val = "Asian" },
    @{ Col = "CG"; Val = "See above under rural section " },
    @{ Col = "CH"; Val = "Native Hawaiian/Pacific Islander" },
    @{ Col = "CI"; Val = "See above under rural section " },
    @{ Col = "CJ"; Val = "Black/African American" },
    @{ Col = "CK"; Val = "See above under rural section " },
    @{ Col = "CL"; Val = "Hispanic/Latino" },
    @{ Col = "CM"; Val = "Our enrollment paperwork is available in Spanish as well as English. See above under rural section for other strategies." },
    @{ Col = "CN"; Val = "Other (list)" },
    @{ Col = "CP"; Val = "Regarding the previous question, what efforts are you making to ensure all of your students (including those from historically underrepresented populations) participate in Outdoor School?" },
    @{ Col = "CQ"; Val = "All students are able to attend Outdoor School through generous donations and school sponsored fundraisers." },
    @{ Col = "CR"; Val = "If you requested and received finding for “extenuating circumstances,” please account for how those funds were used.  Be sure to include documentation." },
    @{ Col = "CS"; Val = "n/a" },
    @{ Col = "CT"; Val = "What local (district, community, partner, parent, etc.) resources or funds have you accessed in addition to your state ODS funding?" },
    @{ Col = "CU"; Val = "In the past, the teacher representing the class that is attending has budjeted for part of Outdoor School with the rest of the funding coming through fundraisers and individual donations." }
)

foreach ($entry in $row4Data) {
    $ws.Range($entry.Col + $rowNum).Value2 = $entry.Val
}

